$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column B
$ws.Range("B5").Value = 187
$ws.Range("B6").Value = 1104
$ws.Range("B7").Value = 1649
$ws.Range("B8").Value = 1200

# Update selection to B8
$ws.Range("B8").Select()
